{"js": "// Remove the \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" footer line, and\n// one of the blank paragraphs that surrounded them (the blank paragraph\n// that used to sit right after the \"LOQ4236: ...\" requisitos line),\n// matching the upstream Jekyll site rebuild that dropped the scraped\n// page-chrome text from the document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the \"LOQ4236: ...\" paragraph; the blank paragraph that\n// immediately follows it (directly before the \"Ver no Jupiter...\" line) is\n// the one that gets collapsed away along with the two text paragraphs.\nlet loqIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4236:\") !== -1) {\n    loqIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (loqIndex !== -1 && items[loqIndex + 1] && items[loqIndex + 1].text === \"\") {\n  toDelete.push(items[loqIndex + 1]);\n}\nfor (const p of items) {\n  if (targetTexts.indexOf(p.text) !== -1) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the scraped page-chrome lines (\"Ver no Jupiter Salvar em pdf\n# Salvar em docx\" and the \"\u00a9 2020 ...\" Jekyll footer) together with the\n# blank paragraph that used to sit right after the \"LOQ4236: ...\" line,\n# matching the latest site rebuild of this course page.\n$d = $word.ActiveDocument\n\n# Locate the \"LOQ4236: ...\" requirement line; the paragraph right after it\n# is the blank paragraph that collapses away along with the two lines below.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"LOQ4236:\")\nif ($found) {\n    $loqPara = $rng.Paragraphs(1)\n    $blankAfter = $loqPara.Next()\n    if ($blankAfter.Range.Text.Trim() -eq \"\") {\n        $blankAfter.Range.Delete()\n    }\n}\n\n# Delete the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph.\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif ($found2) {\n    $rng2.Paragraphs(1).Range.Delete()\n}\n\n# Delete the \"\u00a9 2020 . Contact: ...\" copyright/footer paragraph.\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"Contact: luizeleno@usp.br\")\nif ($found3) {\n    $rng3.Paragraphs(1).Range.Delete()\n}\n"}
